$wb = $excel.ActiveWorkbook

# --- Sheet1: add a "Result" column ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("C1").Value = "Result"
$ws1.Range("C4").Value = "NG"

# --- Sheet2: new sheet with the same layout, different data ---
# Insert it right after Sheet1 so tab order is Sheet1, Sheet2.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "Name"
$ws2.Range("B1").Value = "Code"
$ws2.Range("C1").Value = "Result"

$ws2.Range("A2").Value = "A"
$ws2.Range("B2").Value = 4

$ws2.Range("A3").Value = "B"
$ws2.Range("B3").Value = 5
$ws2.Range("C3").Value = "OK"

$ws2.Range("A4").Value = "C"
$ws2.Range("B4").Value = 6

$ws2.Range("C4").Select() | Out-Null

# Keep Sheet1 as the active/selected sheet as in the source workbook
$ws1.Activate() | Out-Null
$ws1.Select() | Out-Null
$ws1.Range("I2").Select() | Out-Null
